$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.21%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.94%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.103"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.52%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06708"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.16%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.332"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.49%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.447"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.61%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.376"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.30%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9191"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.08%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1587"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.95%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.06800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'6.39%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07710"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.55%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02937"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.02%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.08988"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.14%"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'0.12%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04477"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.79%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0006458"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.59%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.006277"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.51%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.452"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.17%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'2.229"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.55%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.3197"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.49%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1311"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.77%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'4.073"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.73%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D25").Value = "'0.001199"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.25%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004117"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.65%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.10%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001619"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'0.06%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.04264"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'3.78%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006730"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'1.02%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'0.64%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002241"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'3.33%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01193"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.40%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005707"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.13%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.974"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'1.04%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01508"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-18.55%"
$ws.Range("E47").Style = "Normal"
